$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Work from the bottom of the document upward so that earlier paragraph
# indices stay stable while later paragraphs are deleted/merged.
# ---------------------------------------------------------------------

# Paragraph 13: "Link of python code- " + hyperlink -> remove entirely
# (its own Range, including trailing paragraph mark, merges it into the
# trailing empty paragraph 14, which survives as the document's final
# paragraph mark).
$d.Paragraphs.Item(13).Range.Delete()

# Paragraph 12 (empty) -> remove.
$d.Paragraphs.Item(12).Range.Delete()

# Paragraph 11 (empty) -> remove.
$d.Paragraphs.Item(11).Range.Delete()

# Paragraph 10: "Link of the bot-  " + hyperlink -> remove entirely.
$d.Paragraphs.Item(10).Range.Delete()

# Paragraph 9 (empty) -> remove.
$d.Paragraphs.Item(9).Range.Delete()

# Paragraph 8 (empty) -> remove.
$d.Paragraphs.Item(8).Range.Delete()

# Paragraph 7: "- use $inspire to get a nice quote" -> replace with the
# new slider-bar sentence, split into three runs around a spell-check
# marked "and" (<w:proofErr w:type="spellStart"/>and<w:proofErr
# w:type="spellEnd"/>), all inside the SAME paragraph.
$p7 = $d.Paragraphs.Item(7)
$p7sel = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$p7xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">It also has a slider bar which shows the progress of the song and it also shows the total length of the song and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> how much has been convered</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p7sel.InsertXML($p7xml)

# Paragraph 6: "-Whenever you talk about feeling sad or depressed it
# prompts with a supporting message" -> delete just the run text, keep
# the (now empty) paragraph mark.
$p6 = $d.Paragraphs.Item(6)
$p6sel = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$p6sel.Delete()

# Paragraph 5: "-Whenever someone uses a bad work like abuses it checks
# that person with a response message" -> replace text.
$d.Paragraphs.Item(5).Range.Find.Execute("-Whenever someone uses a bad work like abuses it checks that person with a response message", $false, $false, $false, $false, $false, $true, 1, $false, "Can skip, pause, play the previous song, stop.", 2)

# Paragraph 4 (empty) -> remove, merging away (sits between "STARK
# VIGILANCE" and the "Can skip..." paragraph).
$d.Paragraphs.Item(4).Range.Delete()

# Paragraph 3: "STARK VIGILANCE" -> replace text.
$d.Paragraphs.Item(3).Range.Find.Execute("STARK VIGILANCE", $false, $false, $false, $false, $false, $true, 1, $false, "This is a music player which can one song or multiple songs, delete one song or multiple songs.", 2)

# Paragraph 1: "Name of the bot " -> replace text.
$d.Paragraphs.Item(1).Range.Find.Execute("Name of the bot ", $false, $false, $false, $false, $false, $true, 1, $false, "STARK MUSIC PLAYER", 2)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
